# Auto-generated Excel COM-interop edit script
# Commit: "runtime update (2025-10-29 15:08:10)"
# Updates KHL stats workbook: appends 6 new matches to Matches_SOG (rows 401-406,
# dimension A1:G400 -> A1:G406), and refreshes as_of_utc / aggregate stat columns
# on Shots_HA, Shots_Summary and Meta_ext to reflect the 2025-10-28 snapshot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Matches_SOG  -- append new match rows 401-406 and extend dimension
# ---------------------------------------------------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

# Row 401
$c = $wsMatches.Cells.Item(401, 1)
$c.NumberFormat = "@"
$c.Value = "897702"
$c.Style = "Normal"
$wsMatches.Cells.Item(401, 2).Value = "2025-10-28T17:00:00"
$wsMatches.Cells.Item(401, 3).Value = "Трактор"
$wsMatches.Cells.Item(401, 4).Value = "Автомобилист"
$wsMatches.Cells.Item(401, 5).Value = 38
$wsMatches.Cells.Item(401, 6).Value = 25
$wsMatches.Cells.Item(401, 7).Value = "khl_text"

# Row 402
$c = $wsMatches.Cells.Item(402, 1)
$c.NumberFormat = "@"
$c.Value = "897703"
$c.Style = "Normal"
$wsMatches.Cells.Item(402, 2).Value = "2025-10-28T19:00:00"
$wsMatches.Cells.Item(402, 3).Value = "Северсталь"
$wsMatches.Cells.Item(402, 4).Value = "Амур"
$wsMatches.Cells.Item(402, 5).Value = 47
$wsMatches.Cells.Item(402, 6).Value = 20
$wsMatches.Cells.Item(402, 7).Value = "khl_text"

# Row 403
$c = $wsMatches.Cells.Item(403, 1)
$c.NumberFormat = "@"
$c.Value = "897698"
$c.Style = "Normal"
$wsMatches.Cells.Item(403, 2).Value = "2025-10-28T19:10:00"
$wsMatches.Cells.Item(403, 3).Value = "Динамо Мн"
$wsMatches.Cells.Item(403, 4).Value = "ХК Сочи"
$wsMatches.Cells.Item(403, 5).Value = 41
$wsMatches.Cells.Item(403, 6).Value = 18
$wsMatches.Cells.Item(403, 7).Value = "khl_text"

# Row 404
$c = $wsMatches.Cells.Item(404, 1)
$c.NumberFormat = "@"
$c.Value = "897699"
$c.Style = "Normal"
$wsMatches.Cells.Item(404, 2).Value = "2025-10-28T19:30:00"
$wsMatches.Cells.Item(404, 3).Value = "Спартак"
$wsMatches.Cells.Item(404, 4).Value = "ЦСКА"
$wsMatches.Cells.Item(404, 5).Value = 20
$wsMatches.Cells.Item(404, 6).Value = 26
$wsMatches.Cells.Item(404, 7).Value = "khl_text"

# Row 405
$c = $wsMatches.Cells.Item(405, 1)
$c.NumberFormat = "@"
$c.Value = "897700"
$c.Style = "Normal"
$wsMatches.Cells.Item(405, 2).Value = "2025-10-28T19:30:00"
$wsMatches.Cells.Item(405, 3).Value = "Локомотив"
$wsMatches.Cells.Item(405, 4).Value = "Лада"
$wsMatches.Cells.Item(405, 5).Value = 40
$wsMatches.Cells.Item(405, 6).Value = 19
$wsMatches.Cells.Item(405, 7).Value = "khl_text"

# Row 406
$c = $wsMatches.Cells.Item(406, 1)
$c.NumberFormat = "@"
$c.Value = "897701"
$c.Style = "Normal"
$wsMatches.Cells.Item(406, 2).Value = "2025-10-28T19:30:00"
$wsMatches.Cells.Item(406, 3).Value = "Драконы"
$wsMatches.Cells.Item(406, 4).Value = "Динамо М"
$wsMatches.Cells.Item(406, 5).Value = 28
$wsMatches.Cells.Item(406, 6).Value = 22
$wsMatches.Cells.Item(406, 7).Value = "khl_text"

# dimension ref grows from A1:G400 to A1:G406 automatically as Excel recalculates
# the used range once new cell values are written above.

# ---------------------------------------------------------------------------
# Sheet: Shots_HA  -- refresh as_of_utc timestamps and updated shot aggregates
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_HA")

$ws.Range("D2").Value = "2025-10-28T19:30:00Z"
$ws.Range("D3").Value = "2025-10-28T19:30:00Z"
$ws.Range("F3").Value = 24
$ws.Range("K3").Value = 668
$ws.Range("L3").Value = 725
$ws.Range("M3").Value = 27.8
$ws.Range("N3").Value = 30.2
$ws.Range("D4").Value = "2025-10-28T19:30:00Z"
$ws.Range("D5").Value = "2025-10-28T19:30:00Z"
$ws.Range("D6").Value = "2025-10-28T19:30:00Z"
$ws.Range("F6").Value = 20
$ws.Range("K6").Value = 569
$ws.Range("L6").Value = 716
$ws.Range("M6").Value = 28.4
$ws.Range("N6").Value = 35.8
$ws.Range("D7").Value = "2025-10-28T19:30:00Z"
$ws.Range("D8").Value = "2025-10-28T19:30:00Z"
$ws.Range("F8").Value = 19
$ws.Range("K8").Value = 546
$ws.Range("L8").Value = 628
$ws.Range("M8").Value = 28.7
$ws.Range("N8").Value = 33.1
$ws.Range("D9").Value = "2025-10-28T19:30:00Z"
$ws.Range("E9").Value = 19
$ws.Range("G9").Value = 677
$ws.Range("H9").Value = 528
$ws.Range("I9").Value = 35.6
$ws.Range("J9").Value = 27.8
$ws.Range("D10").Value = "2025-10-28T19:30:00Z"
$ws.Range("E10").Value = 19
$ws.Range("G10").Value = 532
$ws.Range("H10").Value = 662
$ws.Range("J10").Value = 34.8
$ws.Range("D11").Value = "2025-10-28T19:30:00Z"
$ws.Range("F11").Value = 15
$ws.Range("K11").Value = 407
$ws.Range("L11").Value = 542
$ws.Range("M11").Value = 27.1
$ws.Range("N11").Value = 36.1
$ws.Range("D12").Value = "2025-10-28T19:30:00Z"
$ws.Range("E12").Value = 15
$ws.Range("G12").Value = 468
$ws.Range("H12").Value = 403
$ws.Range("I12").Value = 31.2
$ws.Range("J12").Value = 26.9
$ws.Range("D13").Value = "2025-10-28T19:30:00Z"
$ws.Range("D14").Value = "2025-10-28T19:30:00Z"
$ws.Range("D15").Value = "2025-10-28T19:30:00Z"
$ws.Range("D16").Value = "2025-10-28T19:30:00Z"
$ws.Range("D17").Value = "2025-10-28T19:30:00Z"
$ws.Range("E17").Value = 14
$ws.Range("G17").Value = 419
$ws.Range("H17").Value = 326
$ws.Range("I17").Value = 29.9
$ws.Range("J17").Value = 23.3
$ws.Range("D18").Value = "2025-10-28T19:30:00Z"
$ws.Range("D19").Value = "2025-10-28T19:30:00Z"
$ws.Range("E19").Value = 24
$ws.Range("G19").Value = 854
$ws.Range("H19").Value = 664
$ws.Range("I19").Value = 35.6
$ws.Range("D20").Value = "2025-10-28T19:30:00Z"
$ws.Range("D21").Value = "2025-10-28T19:30:00Z"
$ws.Range("E21").Value = 15
$ws.Range("G21").Value = 486
$ws.Range("H21").Value = 461
$ws.Range("I21").Value = 32.4
$ws.Range("J21").Value = 30.7
$ws.Range("D22").Value = "2025-10-28T19:30:00Z"
$ws.Range("F22").Value = 18
$ws.Range("K22").Value = 452
$ws.Range("L22").Value = 673
$ws.Range("M22").Value = 25.1
$ws.Range("N22").Value = 37.4
$ws.Range("D23").Value = "2025-10-28T19:30:00Z"
$ws.Range("F23").Value = 19
$ws.Range("K23").Value = 478
$ws.Range("L23").Value = 550
$ws.Range("M23").Value = 25.2
$ws.Range("N23").Value = 28.9

# ---------------------------------------------------------------------------
# Sheet: Shots_Summary  -- refresh as_of_utc timestamps and shot totals
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Shots_Summary")

$ws.Range("D2").Value = "2025-10-28T19:30:00Z"
$ws.Range("D3").Value = "2025-10-28T19:30:00Z"
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 1140
$ws.Range("G3").Value = 1239
$ws.Range("H3").Value = 28.5
$ws.Range("I3").Value = 31
$ws.Range("D4").Value = "2025-10-28T19:30:00Z"
$ws.Range("D5").Value = "2025-10-28T19:30:00Z"
$ws.Range("D6").Value = "2025-10-28T19:30:00Z"
$ws.Range("E6").Value = 36
$ws.Range("F6").Value = 1059
$ws.Range("G6").Value = 1280
$ws.Range("H6").Value = 29.4
$ws.Range("I6").Value = 35.6
$ws.Range("D7").Value = "2025-10-28T19:30:00Z"
$ws.Range("D8").Value = "2025-10-28T19:30:00Z"
$ws.Range("E8").Value = 33
$ws.Range("F8").Value = 1004
$ws.Range("G8").Value = 1024
$ws.Range("H8").Value = 30.4
$ws.Range("I8").Value = 31
$ws.Range("D9").Value = "2025-10-28T19:30:00Z"
$ws.Range("E9").Value = 37
$ws.Range("F9").Value = 1345
$ws.Range("G9").Value = 1012
$ws.Range("H9").Value = 36.4
$ws.Range("I9").Value = 27.4
$ws.Range("D10").Value = "2025-10-28T19:30:00Z"
$ws.Range("E10").Value = 36
$ws.Range("F10").Value = 1012
$ws.Range("G10").Value = 1288
$ws.Range("I10").Value = 35.8
$ws.Range("D11").Value = "2025-10-28T19:30:00Z"
$ws.Range("E11").Value = 37
$ws.Range("F11").Value = 1005
$ws.Range("G11").Value = 1318
$ws.Range("H11").Value = 27.2
$ws.Range("I11").Value = 35.6
$ws.Range("D12").Value = "2025-10-28T19:30:00Z"
$ws.Range("E12").Value = 39
$ws.Range("F12").Value = 1208
$ws.Range("G12").Value = 995
$ws.Range("H12").Value = 31
$ws.Range("I12").Value = 25.5
$ws.Range("D13").Value = "2025-10-28T19:30:00Z"
$ws.Range("D14").Value = "2025-10-28T19:30:00Z"
$ws.Range("D15").Value = "2025-10-28T19:30:00Z"
$ws.Range("D16").Value = "2025-10-28T19:30:00Z"
$ws.Range("D17").Value = "2025-10-28T19:30:00Z"
$ws.Range("E17").Value = 36
$ws.Range("F17").Value = 1141
$ws.Range("G17").Value = 906
$ws.Range("H17").Value = 31.7
$ws.Range("I17").Value = 25.2
$ws.Range("D18").Value = "2025-10-28T19:30:00Z"
$ws.Range("D19").Value = "2025-10-28T19:30:00Z"
$ws.Range("E19").Value = 36
$ws.Range("F19").Value = 1260
$ws.Range("G19").Value = 1112
$ws.Range("H19").Value = 35
$ws.Range("I19").Value = 30.9
$ws.Range("D20").Value = "2025-10-28T19:30:00Z"
$ws.Range("D21").Value = "2025-10-28T19:30:00Z"
$ws.Range("E21").Value = 39
$ws.Range("F21").Value = 1312
$ws.Range("G21").Value = 1245
$ws.Range("H21").Value = 33.6
$ws.Range("I21").Value = 31.9
$ws.Range("D22").Value = "2025-10-28T19:30:00Z"
$ws.Range("E22").Value = 34
$ws.Range("F22").Value = 934
$ws.Range("G22").Value = 1175
$ws.Range("H22").Value = 27.5
$ws.Range("I22").Value = 34.6
$ws.Range("D23").Value = "2025-10-28T19:30:00Z"
$ws.Range("E23").Value = 36
$ws.Range("F23").Value = 861
$ws.Range("G23").Value = 1049
$ws.Range("I23").Value = 29.1

# ---------------------------------------------------------------------------
# Sheet: Meta_ext  -- bump as_of_utc / build_version
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Meta_ext")

$ws.Range("B2").Value = "2025-10-28T19:30:00Z"
$ws.Range("D2").Value = 20

